{"js": "// Fill in the \"In 1990 ...\" forestation placeholder paragraph with the\n// actual queried figures (mirrors the already-filled-in \"In 2016 ...\"\n// paragraph directly above it).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.indexOf(\"In 1990, the percent of the total land area\") === 0) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 1990 forestation placeholder paragraph\");\n}\n\n// Describe the new paragraph as an ordered list of runs: [text, bold, italic, underline, bCs, iCs]\nconst runs = [\n  [\"In 1990, the percent of the total land area of the world designated as forest was \", false, false, false, false, false],\n  [\"3\", true, true, true, true, true],\n  [\"2\", true, true, true, true, true],\n  [\".\", true, true, true, true, true],\n  [\"42\", true, true, true, true, true],\n  [\"%\", true, true, true, true, true],\n  [\".\", true, true, true, true, true],\n  [\" \", true, true, true, true, true],\n  [\"The region with the highest relative forestation was\", false, false, false, false, false],\n  [\" \", true, true, true, false, false],\n  [\"Latin America & Caribbean\", true, true, true, false, false],\n  [\" \", false, false, false, false, false],\n  [\"with \", false, false, false, false, false],\n  [\"51\", true, true, true, false, false],\n  [\".\", true, true, true, false, false],\n  [\"03\", true, true, true, false, false],\n  [\"%\", true, true, true, false, false],\n  [\" \", false, false, false, true, true],\n  [\"and the region with the lowest relative forestation was \", false, false, false, false, false],\n  [\"Middle East & North Africa\", true, true, true, false, false],\n  [\" \", false, false, false, false, false],\n  [\"with \", false, false, false, false, false],\n  [\"1.78\", true, true, true, false, false],\n  [\"%\", true, true, true, false, false],\n  [\" \", false, false, false, true, true],\n  [\"forestation.\", false, false, false, false, false],\n];\n\nfunction escapeXml(s) {\n  return s.replace(/&/g, \"&amp;\").replace(/</g, \"&lt;\").replace(/>/g, \"&gt;\");\n}\n\nlet runsXml = \"\";\nfor (const [text, b, i, u, bCs, iCs] of runs) {\n  let rPr = \"\";\n  if (b || i || u || bCs || iCs) {\n    let props = \"\";\n    if (b) props += \"<w:b/>\";\n    if (bCs) props += \"<w:bCs/>\";\n    if (i) props += \"<w:i/>\";\n    if (iCs) props += \"<w:iCs/>\";\n    if (u) props += '<w:u w:val=\"single\"/>';\n    rPr = `<w:rPr>${props}</w:rPr>`;\n  }\n  const escaped = escapeXml(text);\n  const needsPreserve = text.length === 0 || /^\\s/.test(text) || /\\s$/.test(text);\n  const tEl = needsPreserve ? `<w:t xml:space=\"preserve\">${escaped}</w:t>` : `<w:t>${escaped}</w:t>`;\n  runsXml += `<w:r>${rPr}${tEl}</w:r>`;\n}\n\nconst pkgXml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  `<w:body><w:p>${runsXml}</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\n\n// Clear the paragraph's existing text, then insert the fully-formatted run sequence\n// at its start, so the paragraph's own identity/properties are left untouched.\ntarget.getRange().insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\ntarget.getRange(\"Start\").insertOoxml(pkgXml, Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Fill in the \"In 1990 ...\" forestation placeholder paragraph with the\n# actual queried figures (mirrors the already-filled-in \"In 2016 ...\"\n# paragraph directly above it).\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that still has the 1990 fill-in-the-blank text.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"In 1990, the percent of the total land area*\") {\n        $target = $p\n        break\n    }\n}\n\n# Describe the new paragraph as an ordered list of (text, bold, italic, underline, bCs, iCs) runs.\n$runs = @(\n    @{ Text = \"In 1990, the percent of the total land area of the world designated as forest was \"; B=$false; I=$false; U=$false; BCs=$false; ICs=$false },\n    @{ Text = \"3\";   B=$true; I=$true; U=$true; BCs=$true; ICs=$true },\n    @{ Text = \"2\";   B=$true; I=$true; U=$true; BCs=$true; ICs=$true },\n    @{ Text = \".\";   B=$true; I=$true; U=$true; BCs=$true; ICs=$true },\n    @{ Text = \"42\";  B=$true; I=$true; U=$true; BCs=$true; ICs=$true },\n    @{ Text = \"%\";   B=$true; I=$true; U=$true; BCs=$true; ICs=$true },\n    @{ Text = \".\";   B=$true; I=$true; U=$true; BCs=$true; ICs=$true },\n    @{ Text = \" \";   B=$true; I=$true; U=$true; BCs=$true; ICs=$true },\n    @{ Text = \"The region with the highest relative forestation was\"; B=$false; I=$false; U=$false; BCs=$false; ICs=$false },\n    @{ Text = \" \";   B=$true; I=$true; U=$true; BCs=$false; ICs=$false },\n    @{ Text = \"Latin America & Caribbean\"; B=$true; I=$true; U=$true; BCs=$false; ICs=$false },\n    @{ Text = \" \";   B=$false; I=$false; U=$false; BCs=$false; ICs=$false },\n    @{ Text = \"with \"; B=$false; I=$false; U=$false; BCs=$false; ICs=$false },\n    @{ Text = \"51\";  B=$true; I=$true; U=$true; BCs=$false; ICs=$false },\n    @{ Text = \".\";   B=$true; I=$true; U=$true; BCs=$false; ICs=$false },\n    @{ Text = \"03\";  B=$true; I=$true; U=$true; BCs=$false; ICs=$false },\n    @{ Text = \"%\";   B=$true; I=$true; U=$true; BCs=$false; ICs=$false },\n    @{ Text = \" \";   B=$false; I=$false; U=$false; BCs=$true; ICs=$true },\n    @{ Text = \"and the region with the lowest relative forestation was \"; B=$false; I=$false; U=$false; BCs=$false; ICs=$false },\n    @{ Text = \"Middle East & North Africa\"; B=$true; I=$true; U=$true; BCs=$false; ICs=$false },\n    @{ Text = \" \";   B=$false; I=$false; U=$false; BCs=$false; ICs=$false },\n    @{ Text = \"with \"; B=$false; I=$false; U=$false; BCs=$false; ICs=$false },\n    @{ Text = \"1.78\"; B=$true; I=$true; U=$true; BCs=$false; ICs=$false },\n    @{ Text = \"%\";   B=$true; I=$true; U=$true; BCs=$false; ICs=$false },\n    @{ Text = \" \";   B=$false; I=$false; U=$false; BCs=$true; ICs=$true },\n    @{ Text = \"forestation.\"; B=$false; I=$false; U=$false; BCs=$false; ICs=$false }\n)\n\nfunction Escape-Xml([string]$s) {\n    $s = $s -replace \"&\", \"&amp;\"\n    $s = $s -replace \"<\", \"&lt;\"\n    $s = $s -replace \">\", \"&gt;\"\n    return $s\n}\n\n$runsXml = \"\"\nforeach ($run in $runs) {\n    $rPr = \"\"\n    if ($run.B -or $run.I -or $run.U -or $run.BCs -or $run.ICs) {\n        $props = \"\"\n        if ($run.B)   { $props += \"<w:b/>\" }\n        if ($run.BCs) { $props += \"<w:bCs/>\" }\n        if ($run.I)   { $props += \"<w:i/>\" }\n        if ($run.ICs) { $props += \"<w:iCs/>\" }\n        if ($run.U)   { $props += '<w:u w:val=\"single\"/>' }\n        $rPr = \"<w:rPr>$props</w:rPr>\"\n    }\n    $text = Escape-Xml $run.Text\n    $needsPreserve = ($text.Length -eq 0) -or ($text -match \"^\\s\") -or ($text -match \"\\s$\")\n    if ($needsPreserve) {\n        $runsXml += \"<w:r>$rPr<w:t xml:space=`\"preserve`\">$text</w:t></w:r>\"\n    } else {\n        $runsXml += \"<w:r>$rPr<w:t>$text</w:t></w:r>\"\n    }\n}\n\n$pkgXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Replace the whole paragraph's text span with the new run sequence in one shot.\n$start = $target.Range.Start\n$end = $target.Range.End\n$clearRange = $d.Range($start, $end)\n$clearRange.Text = \"\"\n\n$insertRange = $d.Range($start, $start)\n$insertRange.InsertXML($pkgXml)\n"}
